$d = $word.ActiveDocument
$p18 = $d.Paragraphs(18)
$p19 = $d.Paragraphs(19)

$end19 = $p19.Range.End
$rng = $d.Range($end19 - 1, $end19 - 1)
$rng.InsertParagraphAfter()

$newPara = $d.Paragraphs(20)
$newPara.Range.InsertBefore("Tempat Keberangkatan")

$p20 = $d.Paragraphs(20)
Write-Output ("p20 start=" + $p20.Range.Start + " end=" + $p20.Range.End)
$fullRng = $p20.Range
$fullRng.InsertAfter(" dan Tujuan")
Write-Output ("p20 after: [" + $d.Paragraphs(20).Range.Text + "]")
Write-Output ("p21 after: [" + $d.Paragraphs(21).Range.Text + "]")

$oldRng = $d.Range($p18.Range.Start, $p19.Range.End)
$oldRng.Delete()
Write-Output $d.Paragraphs.Count
